# Add two new lead rows to Sheet1 and refresh the computed Summary sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$summary = $wb.Worksheets.Item("Summary")

# --- Sheet1: new row 4 (Siddharth verma) ---
$ws.Cells.Item(4, 1).Value = "2025-09-14 23:07:08"
$ws.Cells.Item(4, 2).Value = "Siddharth verma"
$ws.Cells.Item(4, 3).Value = "siddharthverma797@gmail.com"
$ws.Cells.Item(4, 4).Value = "graphic Era"
$ws.Cells.Item(4, 5).Value = 112
$ws.Cells.Item(4, 6).Value = "very_hot"
$ws.Cells.Item(4, 7).Value = "/home, /security"
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 330
$ws.Cells.Item(4, 10).Value = 75
$ws.Cells.Item(4, 11).Value = "session_1757871098"

# --- Sheet1: new row 5 (sharaddha email) ---
$ws.Cells.Item(5, 1).Value = "2025-09-14T23:12:28.596886"
$ws.Cells.Item(5, 2).Value = "sharaddha email"
$ws.Cells.Item(5, 3).Value = "sharddha@gmail.com"
$ws.Cells.Item(5, 5).Value = 100
$ws.Cells.Item(5, 6).Value = "Very_hot"
$ws.Cells.Item(5, 7).Value = "/home"
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 0

# --- Sheet1: the old trailing K3 blank cell is gone; row 3 now stops at J3 ---
$ws.Cells.Item(3, 11).ClearContents()

# --- Summary: recomputed aggregate metrics ---
$summary.Cells.Item(2, 2).Value = 3
$summary.Cells.Item(3, 2).Value = "2 (66.7%)"
$summary.Cells.Item(4, 2).Value = 112
$summary.Cells.Item(5, 2).Value = "2025-09-14T17:48:27.117318"
